{"js": "// Office.js (Word JavaScript API) script.\n// Body of: async (context) => { ... }\n//\n// Commit: \"add risco quimico - fumo metalico\"\n//   1) NOME: RENAN NUNES      -> NOME: RENAN FUMOS          (top table)\n//      RENAN NUNES            -> RENAN FUMOS                (signature block)\n//   2) F\u00cdSICO risk list       -> single \" \u2022 N\u00e3o se Aplica\" bullet\n//      QU\u00cdMICO risk list      -> single \" \u2022 Fumos met\u00e1licos\" bullet\n//      ERGON\u00d4MICO risk list   -> single \" \u2022 N\u00e3o se Aplica\" bullet\n//      MEC\u00c2NICO risk list     -> single \" \u2022 N\u00e3o se Aplica\" bullet\n//   3) Parnamirim/RN, 07 de abril de 2025. -> ...09 de abril de 2025.\n\nconst body = context.document.body;\n\n// The four risk-factor cells each hold ONE paragraph whose bullet items\n// are separated by manual line breaks (<w:br/>) inside a single run, not\n// separate paragraphs. `\\u000b` is how Office.js represents that manual\n// line break in paragraph .text. Replacing the whole paragraph's text\n// with \"<newBullet>\\u000b\" keeps that same trailing break and collapses\n// every old bullet/break pair into the one new line.\nasync function replaceBulletParagraph(anchorText, newBullet) {\n  const results = body.search(anchorText, { matchCase: true });\n  results.load(\"text\");\n  await context.sync();\n  const para = results.items[0].paragraphs.getFirst();\n  para.insertText(newBullet + \"\\u000b\", \"Replace\");\n  await context.sync();\n}\n\n// 1) Name: replace both whole-document occurrences of \"RENAN NUNES\".\nconst nameResults = body.search(\"RENAN NUNES\", { matchCase: true });\nnameResults.load(\"text\");\nawait context.sync();\nfor (let i = 0; i < nameResults.items.length; i++) {\n  nameResults.items[i].insertText(\"RENAN FUMOS\", \"Replace\");\n}\nawait context.sync();\n\n// 2) Risk-factor bullet lists.\nawait replaceBulletParagraph(\"Radia\u00e7\u00e3o n\u00e3o ionizante\", \" \u2022 N\u00e3o se Aplica\");\nawait replaceBulletParagraph(\"B\u00e1rio e composto sol\u00faveis como Ba\", \" \u2022 Fumos met\u00e1licos\");\nawait replaceBulletParagraph(\"Esfor\u00e7o f\u00edsico intenso\", \" \u2022 N\u00e3o se Aplica\");\nawait replaceBulletParagraph(\"Animais pe\u00e7onhentos\", \" \u2022 N\u00e3o se Aplica\");\n\n// 3) Footer/closing date line.\nconst dateResults = body.search(\"Parnamirim/RN, 07 de abril de 2025.\", { matchCase: true });\ndateResults.load(\"text\");\nawait context.sync();\ndateResults.items[0].insertText(\"Parnamirim/RN, 09 de abril de 2025.\", \"Replace\");\nawait context.sync();\n", "ps1": "# Word COM interop script: apply \"add risco quimico - fumo metalico\" edit.\n# $word.ActiveDocument is pre-seeded as $d below.\n$d = $word.ActiveDocument\n\n# --- 1) Name fields: \"RENAN NUNES\" -> \"RENAN FUMOS\" (appears twice: once\n#        prefixed with \"NOME: \" in the top table, once alone in the\n#        signature block near the end). A plain Find/Replace on the two\n#        whole-word occurrences handles both. ---\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"RENAN NUNES\"\n$find.Replacement.Text = \"RENAN FUMOS\"\n$find.Forward = $true\n$find.Wrap = 1\n$find.Format = $false\n$find.MatchCase = $true\n$find.MatchWholeWord = $false\n$find.Execute([ref]\"RENAN NUNES\", [ref]$false, [ref]$false, [ref]$false, [ref]$false, [ref]$false, [ref]$true, [ref]1, [ref]$false, \"RENAN FUMOS\", 2) | Out-Null\n\n# --- 2) Risk-factor bullet lists: each is one paragraph inside a table\n#        cell, where the bullet items are separated by manual line breaks\n#        (vertical-tab / Chr(11)) inside a single run, not real paragraph\n#        marks. Replace the whole bullet-list paragraph's content (but\n#        keep its own trailing line break) with a single new bullet. ---\nfunction Set-BulletParagraph($doc, [string]$anchorText, [string]$newBullet) {\n    foreach ($p in $doc.Paragraphs) {\n        $t = $p.Range.Text\n        if ($t -ne $null -and $t.Contains($anchorText)) {\n            $pr = $p.Range\n            # Range.End counts the trailing paragraph/cell mark as a single\n            # unit (even though it renders as \\r + cell-mark) -- trim just\n            # one unit so the paragraph's own last manual line break stays.\n            $newEnd = $pr.End - 1\n            $sub = $doc.Range($pr.Start, $newEnd)\n            $sub.Text = $newBullet + [char]11\n            return $true\n        }\n    }\n    return $false\n}\n\nSet-BulletParagraph $d \"Radia\u00e7\u00e3o n\u00e3o ionizante\" \" \u2022 N\u00e3o se Aplica\" | Out-Null\nSet-BulletParagraph $d \"B\u00e1rio e composto sol\u00faveis como Ba\" \" \u2022 Fumos met\u00e1licos\" | Out-Null\nSet-BulletParagraph $d \"Esfor\u00e7o f\u00edsico intenso\" \" \u2022 N\u00e3o se Aplica\" | Out-Null\nSet-BulletParagraph $d \"Animais pe\u00e7onhentos\" \" \u2022 N\u00e3o se Aplica\" | Out-Null\n\n# --- 3) Footer date: \"07 de abril de 2025\" -> \"09 de abril de 2025\" ---\n$find2 = $d.Content.Find\n$find2.ClearFormatting()\n$find2.Replacement.ClearFormatting()\n$find2.Text = \"Parnamirim/RN, 07 de abril de 2025.\"\n$find2.Replacement.Text = \"Parnamirim/RN, 09 de abril de 2025.\"\n$find2.Forward = $true\n$find2.Wrap = 1\n$find2.Format = $false\n$find2.MatchCase = $true\n$find2.Execute([ref]\"Parnamirim/RN, 07 de abril de 2025.\", [ref]$false, [ref]$false, [ref]$false, [ref]$false, [ref]$false, [ref]$true, [ref]1, [ref]$false, \"Parnamirim/RN, 09 de abril de 2025.\", 2) | Out-Null\n"}
